# edit.ps1
# Applies the "TwelveBarBreakout" results update:
#   - "Signal Instances" sheet: append 12 new rows (34-45) of
#     twelve_bar_breakout signal instances.
#   - "Summary" sheet: insert a new summary row for twelve_bar_breakout
#     above the existing random_baseline row (which shifts from row 4
#     to row 5).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: write an array of values into a worksheet row starting at
# column A. $null entries are left untouched (i.e. stay blank), which
# mirrors the source data's blank / inline-empty-string cells.
# ---------------------------------------------------------------------
function Set-RowValues {
    param($ws, [int]$row, [object[]]$values)
    for ($i = 0; $i -lt $values.Count; $i++) {
        $value = $values[$i]
        if ($null -ne $value) {
            $col = $i + 1
            $ws.Cells.Item($row, $col).Value = $value
        }
    }
}

# ---------------------------------------------------------------------
# 1. "Signal Instances" sheet (sheet index 1) - append rows 34-45.
#    Columns: signal_type, signal_bar, signal_date, entry_price,
#    return_12m, profitable_12m, mfe_12m, mfe_bar, mae_12m, mae_bar,
#    exit_signal_fired, exit_signal_bar, return_at_exit, left_on_table,
#    exit_vs_hold, exit_vs_mfe
# ---------------------------------------------------------------------
$wsSignals = $wb.Worksheets.Item("Signal Instances")

$sheet1Rows = @{}
$sheet1Rows[34] = @("twelve_bar_breakout", 29, $null, 10.20273590087891, 28.10489174773218, $true, 45.65791965970256, 40, -1.844731544352345, 1, $true, 29, 18.63229486731778, 17.55302791197039, -9.472596880414393, -27.02562479238478)
$sheet1Rows[35] = @("twelve_bar_breakout", 61, $null, 13.12647724151611, 15.32812621579793, $true, 20.53698128851986, 29, -13.00611971586399, 5, $false, $null, $null, 5.208855072721928, $null, $null)
$sheet1Rows[36] = @("twelve_bar_breakout", 86, $null, 14.9103832244873, -2.670916623120283, $false, 6.115712642873021, 4, -14.34800560110058, 32, $true, 31, -8.253970182933186, 8.786629265993305, -5.583053559812902, -14.36968282580621)
$sheet1Rows[37] = @("twelve_bar_breakout", 111, $null, 15.68036937713623, -4.386287629791378, $false, 2.808901465063167, 33, -18.55395560442638, 7, $true, 6, -12.75916842288426, 7.195189094854546, -8.372880793092882, -15.56806988794743)
$sheet1Rows[38] = @("twelve_bar_breakout", 139, $null, 14.98878765106201, 10.87225119637272, $true, 12.34192574049123, 52, -9.608032836994637, 27, $true, 14, -1.128366228850848, 1.469674544118517, -12.00061742522356, -13.47029196934208)
$sheet1Rows[39] = @("twelve_bar_breakout", 178, $null, 16.59197235107422, 0.3767111792060027, $true, 4.354744652040772, 17, -13.73212299804742, 28, $true, 21, -10.32894461552469, 3.978033472834769, -10.70565579473069, -14.68368926756546)
$sheet1Rows[40] = @("twelve_bar_breakout", 218, $null, 15.68484401702881, -7.446911551102099, $false, 9.913803646290392, 11, -11.96169864274231, 39, $true, 18, 3.907890379078518, 17.36071519739249, 11.35480193018062, -6.005913267211874)
$sheet1Rows[41] = @("twelve_bar_breakout", 228, $null, 16.80080795288086, -7.139401861528164, $false, 2.612973753728526, 1, -17.80948700886833, 29, $true, 8, -2.994007353383599, 9.752375615256689, 4.145394508144565, -5.606981107112125)
$sheet1Rows[42] = @("twelve_bar_breakout", 273, $null, 15.02927398681641, 26.97776099196058, $true, 27.49315692572213, 52, -4.835691816632501, 1, $true, 37, 13.44894226508021, 0.5153959337615568, -13.52881872688037, -14.04421466064193)
$sheet1Rows[43] = @("twelve_bar_breakout", 290, $null, 16.17246055603027, 35.56898686606376, $true, 37.46036945682945, 48, -1.353734816567277, 0, $true, 20, 5.429549876418603, 1.891382590765694, -30.13943698964516, -32.03081958041085)
$sheet1Rows[44] = @("twelve_bar_breakout", 294, $null, 16.53417587280273, 22.09994151351938, $true, 34.45317263829772, 44, -2.36040146776894, 0, $true, 16, 3.123085779018004, 12.35323112477834, -18.97685573450137, -31.33008685927971)
$sheet1Rows[45] = @("twelve_bar_breakout", 300, $null, 17.67643356323242, 16.38682619404367, $true, 25.76475877362436, 38, -5.293573323669627, 11, $true, 10, -3.540766257115425, 9.37793257958069, -19.92759245115909, -29.30552503073978)

foreach ($rowNum in ($sheet1Rows.Keys | Sort-Object)) {
    Set-RowValues $wsSignals $rowNum $sheet1Rows[$rowNum]
}

# ---------------------------------------------------------------------
# 2. "Summary" sheet (sheet index 2) - insert a new row 4 for
#    twelve_bar_breakout, pushing the existing random_baseline row
#    (previously row 4) down to row 5.
#    Columns: signal_type, total_signals, win_rate_12m, mean_return_12m,
#    median_return_12m, std_return_12m, mean_mfe_12m, mean_mae_12m,
#    mean_left_on_table, exit_fired_rate, mean_exit_bar,
#    exit_useful_rate, mean_exit_vs_hold, total_samples,
#    baseline_win_rate, baseline_mean_return, baseline_median_return,
#    baseline_std_return, baseline_mean_mfe, baseline_mean_mae
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Rows.Item(4).Insert()

$sheet2Row4 = @("twelve_bar_breakout", 12, 66.66666666666666, 11.17266485326286, 13.10018870608532, 15.15955036377552, 19.12620172026527, -9.558962948086196, 7.95353686700241, 91.66666666666666, 19.09090909090909, 18.18181818181818, -10.29157381064866, $null, $null, $null, $null, $null, $null, $null)

Set-RowValues $wsSummary 4 $sheet2Row4
